# Each data row (1-6) has a set of readout lines in columns C:J
# (Concentration / particle-size counts / AT-RH / DP-WB). The author
# changed these from single-line strings to "looking" (line-broken)
# text: every line gets a trailing line break so the cell displays as
# its own paragraph, and the two instrument-reading lines (the
# "AT: ... RH: ..." and "DP: ...WB: ..." lines) additionally keep a
# couple of trailing spaces before that line break.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 1; $row -le 6; $row++) {
  for ($col = 3; $col -le 10; $col++) {
    $cell = $ws.Cells.Item($row, $col)
    $text = [string]$cell.Value2

    if ($text.StartsWith("AT:") -or $text.StartsWith("DP:")) {
      $cell.Value = $text + "  `n"
    } else {
      $cell.Value = $text + "`n"
    }
  }
}
